$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price values stay as text (matching source data),
# same as the original cells which are stored as text.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '28.220.46'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '1.877.24'
$ws.Range('E3').Value = '  +1.69%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.51%  '
$ws.Range('D5').Value = '316.11'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +1.71%  '
$ws.Range('D8').Value = '0.3714'
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('D9').Value = '0.07422'
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('D10').Value = '0.8841'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('E11').Value = '  +2.58%  '
$ws.Range('D12').Value = '1.918.82'
$ws.Range('E12').Value = '  +3.27%  '
$ws.Range('D13').Value = '5.501'
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('D14').Value = '6.637'
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '0.06984'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '81.38'
$ws.Range('E17').Value = '  +3.45%  '
$ws.Range('E18').Value = '  +3.46%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = '15.60'
$ws.Range('E20').Value = '  +1.70%  '
$ws.Range('D21').Value = '28.311.43'
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').Value = '5.081'
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').Value = '10.90'
$ws.Range('E23').Value = '  +5.24%  '
$ws.Range('D24').Value = '2.145.83'
$ws.Range('E24').Value = '  +4.04%  '
$ws.Range('D25').Value = '1.969'
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('D26').Value = '154.02'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').Value = '18.80'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').Value = '5.441'
$ws.Range('E28').Value = '  +3.73%  '
$ws.Range('D29').Value = '117.77'
$ws.Range('E29').Value = '  -2.86%  '
$ws.Range('D30').Value = '1.872'
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').Value = '0.08992'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').Value = '0.7975'
$ws.Range('E32').Value = '  +5.18%  '
$ws.Range('D33').Value = '4.731'
$ws.Range('E33').Value = '  +4.16%  '
$ws.Range('D34').Value = '1.186'
$ws.Range('E34').Value = '  +8.34%  '
$ws.Range('D35').Value = '2.934'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.129'
$ws.Range('E37').Value = '  +3.36%  '
$ws.Range('D38').Value = '0.05471'
$ws.Range('E38').Value = '  +2.27%  '
$ws.Range('D39').Value = '0.01960'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').Value = '2.882'
$ws.Range('E40').Value = '  +2.67%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '0.1695'
$ws.Range('E41').Value = '  +3.04%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.5178'
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('D43').Value = '6.907'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '8.634'
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('E45').Value = '  +2.02%  '
$ws.Range('D46').Value = '0.4783'
$ws.Range('E46').Value = '  +0.95%  '
$ws.Range('D47').Value = '0.06598'
$ws.Range('E47').Value = '  +0.91%  '
$ws.Range('D48').Value = '105.93'
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('D49').Value = '1.001'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('D50').Value = '1.657'
$ws.Range('E50').Value = '  +2.36%  '
$ws.Range('D51').Value = '1.861'
$ws.Range('E51').Value = '  +8.17%  '
